$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.775.90'
$ws.Cells.Item(2, 5).Value = '  +0.29%  '
$ws.Cells.Item(3, 4).Value = '1.636.18'
$ws.Cells.Item(3, 5).Value = '  +0.38%  '
$ws.Cells.Item(4, 5).Value = '  -0.66%  '
$ws.Cells.Item(5, 4).Formula = '''216.61'
$ws.Cells.Item(5, 5).Value = '  -1.05%  '
$ws.Cells.Item(6, 4).Formula = '''0.504'
$ws.Cells.Item(6, 5).Value = '  +1.88%  '
$ws.Cells.Item(7, 5).Value = '  -0.66%  '
$ws.Cells.Item(8, 4).Formula = '''0.253'
$ws.Cells.Item(8, 5).Value = '  +1.98%  '
$ws.Cells.Item(9, 4).Formula = '''0.0622'
$ws.Cells.Item(9, 5).Value = '  +0.79%  '
$ws.Cells.Item(10, 4).Formula = '''19.84'
$ws.Cells.Item(10, 5).Value = '  +5.32%  '
$ws.Cells.Item(11, 4).Formula = '''0.0844'
$ws.Cells.Item(11, 5).Value = '  +0.00%  '
$ws.Cells.Item(12, 4).Value = '1.865.28'
$ws.Cells.Item(12, 5).Value = '  +0.35%  '
$ws.Cells.Item(13, 4).Value = '1.652.42'
$ws.Cells.Item(13, 5).Value = '  +0.51%  '
$ws.Cells.Item(14, 5).Value = '  +0.45%  '
$ws.Cells.Item(15, 5).Value = '  +1.86%  '
$ws.Cells.Item(16, 4).Formula = '''66.30'
$ws.Cells.Item(16, 5).Value = '  +3.60%  '
$ws.Cells.Item(17, 4).Value = '26.779.76'
$ws.Cells.Item(17, 5).Value = '  +0.36%  '
$ws.Cells.Item(18, 4).Value = '0.0₃0728'
$ws.Cells.Item(19, 4).Formula = '''218.17'
$ws.Cells.Item(19, 5).Value = '  +2.73%  '
$ws.Cells.Item(20, 5).Value = '  -0.64%  '
$ws.Cells.Item(21, 2).Value = 'Chainlink'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(21, 4).Formula = '''6.62'
$ws.Cells.Item(21, 5).Value = '  +6.85%  '
$ws.Cells.Item(22, 2).Value = 'Uniswap'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(22, 4).Formula = '''4.37'
$ws.Cells.Item(22, 5).Value = '  +1.64%  '
$ws.Cells.Item(23, 4).Formula = '''2.42'
$ws.Cells.Item(23, 5).Value = '  +3.74%  '
$ws.Cells.Item(24, 4).Formula = '''9.12'
$ws.Cells.Item(24, 5).Value = '  +1.89%  '
$ws.Cells.Item(25, 4).Formula = '''146.92'
$ws.Cells.Item(25, 5).Value = '  -0.31%  '
$ws.Cells.Item(26, 5).Value = '  -0.63%  '
$ws.Cells.Item(27, 5).Value = '  +5.99%  '
$ws.Cells.Item(28, 5).Value = '  +1.23%  '
$ws.Cells.Item(29, 4).Formula = '''15.74'
$ws.Cells.Item(29, 5).Value = '  +1.16%  '
$ws.Cells.Item(30, 5).Value = '  +0.99%  '
$ws.Cells.Item(31, 5).Value = '  -0.99%  '
$ws.Cells.Item(32, 4).Formula = '''3.33'
$ws.Cells.Item(32, 5).Value = '  -0.79%  '
$ws.Cells.Item(33, 5).Value = '  +1.29%  '
$ws.Cells.Item(34, 4).Formula = '''1.55'
$ws.Cells.Item(34, 5).Value = '  +1.95%  '
$ws.Cells.Item(35, 5).Value = '  -0.25%  '
$ws.Cells.Item(36, 4).Value = '1.245.29'
$ws.Cells.Item(36, 5).Value = '  -0.31%  '
$ws.Cells.Item(37, 4).Formula = '''0.0176'
$ws.Cells.Item(37, 5).Value = '  +1.14%  '
$ws.Cells.Item(38, 5).Value = '  +2.23%  '
$ws.Cells.Item(39, 4).Formula = '''0.829'
$ws.Cells.Item(39, 5).Value = '  +3.81%  '
$ws.Cells.Item(40, 5).Value = '  -0.61%  '
$ws.Cells.Item(41, 4).Formula = '''0.807'
$ws.Cells.Item(41, 5).Value = '  +0.81%  '
$ws.Cells.Item(42, 4).Formula = '''5.37'
$ws.Cells.Item(42, 5).Value = '  +2.49%  '
$ws.Cells.Item(43, 4).Value = '1.777.24'
$ws.Cells.Item(43, 5).Value = '  +0.50%  '
$ws.Cells.Item(44, 4).Formula = '''61.73'
$ws.Cells.Item(44, 5).Value = '  +3.94%  '
$ws.Cells.Item(45, 4).Formula = '''2.09'
$ws.Cells.Item(45, 5).Value = '  -1.72%  '
$ws.Cells.Item(46, 4).Formula = '''91.42'
$ws.Cells.Item(46, 5).Value = '  -0.11%  '
$ws.Cells.Item(47, 5).Value = '  +0.68%  '
$ws.Cells.Item(48, 5).Value = '  +4.14%  '
$ws.Cells.Item(49, 5).Value = '  -0.38%  '
$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).Formula = '''7.63'
$ws.Cells.Item(50, 5).Value = '  +2.95%  '
$ws.Cells.Item(51, 2).Value = 'Algorand'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(51, 4).Formula = '''0.0968'
$ws.Cells.Item(51, 5).Value = '  +1.95%  '
